$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF") - reuse H1 formatting (bold, centered, bordered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# I0 and IF values for rows 2-54
$values = @{
    2 = @(9, 9)
    3 = @(9, 9)
    4 = @(8, 8)
    5 = @(9, 9)
    6 = @(9, 9)
    7 = @(9, 9)
    8 = @(9, 9)
    9 = @(8, 8)
    10 = @(9, 9)
    11 = @(9, 9)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(9, 9)
    16 = @(9, 9)
    17 = @(7, 7)
    18 = @(6, 6)
    19 = @(9, 9)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(8, 8)
    23 = @(9, 9)
    24 = @(7, 7)
    25 = @(5, 6)
    26 = @(8, 8)
    27 = @(10, 11)
    28 = @(7, 7)
    29 = @(5, 6)
    30 = @(12, 12)
    31 = @(6, 7)
    32 = @(8, 8)
    33 = @(6, 6)
    34 = @(7, 7)
    35 = @(8, 8)
    36 = @(8, 8)
    37 = @(5, 5)
    38 = @(8, 8)
    39 = @(7, 7)
    40 = @(8, 8)
    41 = @(7, 7)
    42 = @(9, 9)
    43 = @(6, 6)
    44 = @(9, 9)
    45 = @(8, 8)
    46 = @(9, 9)
    47 = @(8, 8)
    48 = @(6, 6)
    49 = @(8, 8)
    50 = @(4, 4)
    51 = @(7, 7)
    52 = @(4, 4)
    53 = @(6, 6)
    54 = @(4, 5)
}

foreach ($r in $values.Keys) {
    $pair = $values[$r]
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
}
